$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new working-hours entry (2014-07-12, 14:45-18:00) was recorded. It lands
# right after the existing row 133 entry, so insert a fresh row at 134 which
# shifts the trailing "totals" rows (134-137) down to (135-138) and adjusts
# all of the formulas/references that live in those rows automatically.
$ws.Rows("134:134").Insert()

# The previous last entry (row 133) actually ran until 13:00 instead of noon.
$ws.Range("E133").Value = 0.54166666666666663

# Fill in the new entry in row 134.
$ws.Range("A134").Value = 2014
$ws.Range("B134").Value = 7
$ws.Range("C134").Value = 12
$ws.Range("D134").Value = 0.61458333333333337
$ws.Range("E134").Value = 0.75
$ws.Range("F134").Formula = "=(E134-D134)*24*60"
$ws.Range("G134").Formula = "=F134/60"

# The "sum [min]" total (now on row 136) needs to include the new row.
$ws.Range("F136").Formula = "=SUM(F2:F134)"

# Update the selection/scroll position to reflect where the user ended up
# after the edit.
$ws.Range("E135").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 118
$win.ScrollColumn = 1
